$d = $word.ActiveDocument

# --- Locate the relevant skills paragraphs by their current text -----------
$basesIndex = -1
$mlopsIndex = -1
$visuIndex = -1

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Bases de donn*es : SQL, MongoDB, Neo4j, Redis*") {
        $basesIndex = $p.Index
    } elseif ($t -like "MLOps : aws, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit*") {
        $mlopsIndex = $p.Index
    } elseif ($t -like "Visualisation : tableau*") {
        $visuIndex = $p.Index
    }
}

# --- 1) Insert a new (cloned-formatting) paragraph right after the
#        "Bases de données" line, then populate it with the MLOps text. ----
$basesPara = $d.Paragraphs.Item($basesIndex)
$basesPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($basesIndex + 1)
$newPara.Range.Text = "MLOps : aws, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit"

# --- 2) Turn the original "Bases de données" paragraph into "Visualisation". 
$basesPara = $d.Paragraphs.Item($basesIndex)
$basesPara.Range.Text = "Visualisation : tableau"

# --- 3) The old "Visualisation : tableau" paragraph shifted down by one
#        slot because of the insertion above; delete that duplicate. -------
$visuPara = $d.Paragraphs.Item($visuIndex + 1)
$visuPara.Range.Delete()

# --- 4) Turn the old "MLOps" paragraph into "Bases de données" (its index
#        is unaffected since the insertion/deletion happened before it). ---
$mlopsPara = $d.Paragraphs.Item($mlopsIndex)
$mlopsPara.Range.Text = "Bases de données : SQL, MongoDB, Neo4j, Redis"
